$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.883.26"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.877.21"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'324.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.4608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").Value = "'0.3863"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("D9").Value = "'0.07856"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").Value = "'0.9843"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "'21.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").Value = "1.890.64"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "'6.982"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").Value = "'5.645"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "'0.06970"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "'87.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "'0.000009944"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "'16.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "28.894.72"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").Value = "'10.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "'2.105"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("D25").Value = "'156.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "'19.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("D27").Value = "'5.971"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'117.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "'1.908"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'0.09356"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").Value = "'0.9006"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("D32").Value = "'5.258"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").Value = "'1.317"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("D34").Value = "'3.252"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Value = "'0.02072"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("D40").Value = "'0.5643"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("D41").Value = "'0.1763"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").Value = "'9.691"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").Value = "'2.251"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "'11.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'0.5333"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").Value = "'0.07039"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("D47").Value = "'1.838"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").Value = "'2.565"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D49").Value = "'112.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "'1.063"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.86%  "
$ws.Range("D51").Value = "'70.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.68%  "
